$d = $word.ActiveDocument
$tab = [char]9

# ---------------------------------------------------------------------------
# Locate paragraphs by their distinctive text rather than hard-coded indices,
# since indices shift as we insert/delete content.
# ---------------------------------------------------------------------------
function Find-ParagraphIndexByStart($doc, [string]$needle) {
    $n = $doc.Paragraphs.Count
    for ($i = 1; $i -le $n; $i++) {
        $t = $doc.Paragraphs.Item($i).Range.Text
        if ($t.Contains($needle)) {
            return $i
        }
    }
    return -1
}

$goalsIdx    = Find-ParagraphIndexByStart $d "We have three main goals here"
$thousandIdx = Find-ParagraphIndexByStart $d "If the girl counts from 1 to 1000"
$iUsedIdx    = Find-ParagraphIndexByStart $d "I used "

# ---------------------------------------------------------------------------
# STEP 1: Split the "I used ..." paragraph right after its leading tab
# character (i.e. right before the _GoBack bookmark), so the tab stays
# behind in its own paragraph and "I used ..." (with the bookmark) becomes
# its own paragraph.
# ---------------------------------------------------------------------------
$iUsedPara = $d.Paragraphs.Item($iUsedIdx)
$splitPoint = $d.Range($iUsedPara.Range.Start + 1, $iUsedPara.Range.Start + 1)
$splitPoint.InsertParagraphAfter()

# The paragraph at $iUsedIdx is now just the lone tab; the paragraph right
# after it holds the bookmark + "I used ..." text.
$tabOnlyPara = $d.Paragraphs.Item($iUsedIdx)

# ---------------------------------------------------------------------------
# STEP 2: Append the new "sub goal" sentence onto the lone-tab paragraph.
# (Use End-1 so we land *inside* the paragraph, before its paragraph mark.)
# ---------------------------------------------------------------------------
$appendPoint = $d.Range($tabOnlyPara.Range.End - 1, $tabOnlyPara.Range.End - 1)
$appendPoint.InsertAfter("The sub goal is to know in witch finger she stops.")

# ---------------------------------------------------------------------------
# STEP 3 + 4: Insert a blank paragraph, then the re-created
# "goals / also / 100 / 1000 / blank" block, right after the new "sub goal"
# paragraph (still before "I used ...").
# ---------------------------------------------------------------------------
$subGoalPara = $d.Paragraphs.Item($iUsedIdx)
$insertPt = $d.Range($subGoalPara.Range.End - 1, $subGoalPara.Range.End - 1)

$goalsText    = "We have three main goals here: based on this sequence of instructions which finger will be identified with the ending count of 10, 100 and 1000. I believe we need to build a sequential code of the five-finger counts. Then we would implement a conditional code to change up the sequence with a set of instructions to reverse the order of count. Then back to the sequential code for a count of five. This program will last until the first total equals 10. Then we change the next program to hit a total 100. While the third program would hit a total of 1000."
$alsoText     = "Also, if the girl counts from one to ten, she will stop at the first finger."
$hundredText  = "If the girl counts from 1 to 100, she will stop at the ring finger."
$thousandText = "If the girl counts from 1 to 1000, she will stop at the first finger."

$block = "`r${tab}${goalsText}`r`r${tab}${alsoText}`r${tab}${hundredText}`r${tab}${thousandText}`r"
$insertPt.InsertAfter($block)

# ---------------------------------------------------------------------------
# STEP 5: Delete the original "goals / also / 100 / 1000 / blank" block that
# is still sitting above (right before the lone-tab / sub-goal paragraph),
# including the blank paragraph that originally separated it from "I used ...".
# ---------------------------------------------------------------------------
$goalsPara          = $d.Paragraphs.Item($goalsIdx)
$emptyAfterThousand = $d.Paragraphs.Item($thousandIdx + 1)
$delRange = $d.Range($goalsPara.Range.Start, $emptyAfterThousand.Range.End)
$delRange.Delete()

# ---------------------------------------------------------------------------
# STEP 6: Append the closing "Yes each solution ..." paragraphs at the very
# end of the document (after "... first finger.").
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$closingText = "Yes each solution meets the goal.Yes each solution work for all cases because it will always stop at the same finger."
$endPt = $d.Range($lastPara.Range.End - 1, $lastPara.Range.End - 1)
$endPt.InsertAfter("`r`r${tab}${closingText}")

Write-Host "Final paragraph count:" $d.Paragraphs.Count
